# Apply updated crypto price/volume data per GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting (prices/links/percent strings must stay text, not be
# auto-converted to numbers/dates by Excel) by forcing the Text number format
# on each cell before assigning its new value.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.976.55"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.294.22"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.24"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.99"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.59%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.96%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.509"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.71%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +8.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0788"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.09%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.19"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +7.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.93"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.654.17"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.799"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.85%  "
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.205.30"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -4.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.878.99"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.50"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +8.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0901"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.09"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.15%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.63"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +10.69%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.43"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.94"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.72%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +15.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.38"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "166.94"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.00"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.82%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +4.39%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.36%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0685"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.75%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.35%  "
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.78"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.16%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.81"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.16%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "ApeXProtocol"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.29"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.29%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0291"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.964.73"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.41%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.50%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.50"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.35"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +4.69%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.54"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.93%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.519.92"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "70.65"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.22%  "
